$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "SCD0011"
$ws.Range("B2").Value = "SCD0011-028"
$ws.Columns.Item(2).ColumnWidth = 11.7
$ws.Range("B3").Select()
